# Update "tasas-transfi.xlsx" with latest automated rate values.

$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 2.6 = 9623.38 pesos
✅ 9623.38 pesos = 2.58 = 927.51 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures in N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 385
$wsTasas.Range("O10").Value = 3705
$wsTasas.Range("N12").Value = 3730
$wsTasas.Range("O12").Value = 359.5
